$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.210.29"
$ws.Range("E2").Value = "  -2.94%  "
$ws.Range("D3").Value = "'1.916.67"
$ws.Range("E3").Value = "  -3.55%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -1.18%  "
$ws.Range("D5").Value = "'327.91"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "'0.4683"
$ws.Range("E7").Value = "  -5.88%  "
$ws.Range("D8").Value = "'0.4023"
$ws.Range("E8").Value = "  -4.11%  "
$ws.Range("D9").Value = "'53.22"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").Value = "'0.08408"
$ws.Range("E10").Value = "  -9.77%  "
$ws.Range("D11").Value = "'1.046"
$ws.Range("E11").Value = "  -4.46%  "
$ws.Range("D12").Value = "'22.14"
$ws.Range("E12").Value = "  -4.39%  "
$ws.Range("D13").Value = "'1.918.62"
$ws.Range("E13").Value = "  -4.56%  "
$ws.Range("D14").Value = "'7.414"
$ws.Range("E14").Value = "  -7.06%  "
$ws.Range("D15").Value = "'6.074"
$ws.Range("E15").Value = "  -5.72%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "'89.68"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("D18").Value = "'0.00001060"
$ws.Range("E18").Value = "  -4.71%  "
$ws.Range("D19").Value = "'0.06596"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").Value = "'17.97"
$ws.Range("E20").Value = "  -8.11%  "
$ws.Range("D21").Value = "'0.9993"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").Value = "'5.719"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").Value = "'28.199.86"
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("D24").Value = "'11.27"
$ws.Range("E24").Value = "  -5.95%  "
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "'2.130.31"
$ws.Range("E26").Value = "  -5.15%  "
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "'20.00"
$ws.Range("E28").Value = "  -3.88%  "
$ws.Range("D29").Value = "'5.737"
$ws.Range("E29").Value = "  -8.49%  "
$ws.Range("E30").Value = "  -6.23%  "
$ws.Range("D31").Value = "'123.49"
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("D32").Value = "'0.9747"
$ws.Range("E32").Value = "  -6.93%  "
$ws.Range("D33").Value = "'0.09619"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").Value = "'1.452"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("D35").Value = "'3.633"
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("D36").Value = "'5.546"
$ws.Range("E36").Value = "  -4.57%  "
$ws.Range("D37").Value = "'0.02302"
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.264"
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'8.770"
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").Value = "'0.6137"
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("D42").Value = "'11.05"
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").Value = "'0.1905"
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("D45").Value = "'1.307"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5848"
$ws.Range("E46").Value = "  -5.98%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.74"
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("D48").Value = "'2.025"
$ws.Range("E48").Value = "  -7.25%  "
$ws.Range("D49").Value = "'3.438"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").Value = "'0.06887"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("E51").Value = "  -9.16%  "
